# Update "Horarios Línea 141" workbook with the latest scrape (765).
# Touches all three sheets: LP1912, LP1912-215 and 6203-6173.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Cells.Item(2, 1).Value = "Última actualización: 03:21:41"
$ws1.Cells.Item(3, 1).Value = "Total filas: 16"

$ws1.Cells.Item(12, 1).Value = "03:21:41"
$ws1.Cells.Item(12, 2).Value = "03:24"
$ws1.Cells.Item(12, 3).Value = "14_ABASTO"
$ws1.Cells.Item(12, 4).Value = 3
$ws1.Cells.Item(12, 5).Value = "LP1912"

$ws1.Cells.Item(13, 1).Value = "02:55:01"
$ws1.Cells.Item(13, 2).Value = "03:48"
$ws1.Cells.Item(13, 3).Value = "14_ABASTO"
$ws1.Cells.Item(13, 4).Value = 53
$ws1.Cells.Item(13, 5).Value = "LP1912"

$ws1.Cells.Item(14, 1).Value = "01:59:40"
$ws1.Cells.Item(14, 2).Value = "03:50"
$ws1.Cells.Item(14, 3).Value = "14_ABASTO"
$ws1.Cells.Item(14, 4).Value = 111
$ws1.Cells.Item(14, 5).Value = "LP1912"

$ws1.Cells.Item(15, 1).Value = "02:30:53"
$ws1.Cells.Item(15, 2).Value = "03:52"
$ws1.Cells.Item(15, 3).Value = "14_ABASTO"
$ws1.Cells.Item(15, 4).Value = 82
$ws1.Cells.Item(15, 5).Value = "LP1912"

$ws1.Cells.Item(16, 1).Value = "03:21:41"
$ws1.Cells.Item(16, 2).Value = "04:01"
$ws1.Cells.Item(16, 3).Value = "81_EL PELIGRO"
$ws1.Cells.Item(16, 4).Value = 40
$ws1.Cells.Item(16, 5).Value = "LP1912"

$ws1.Cells.Item(17, 1).Value = "03:21:41"
$ws1.Cells.Item(17, 2).Value = "04:45"
$ws1.Cells.Item(17, 3).Value = "215A_EL PATO"
$ws1.Cells.Item(17, 4).Value = 84
$ws1.Cells.Item(17, 5).Value = "LP1912"

$ws1.Cells.Item(18, 1).Value = "02:55:01"
$ws1.Cells.Item(18, 2).Value = "04:46"
$ws1.Cells.Item(18, 3).Value = "215A_EL PATO"
$ws1.Cells.Item(18, 4).Value = 111
$ws1.Cells.Item(18, 5).Value = "LP1912"

$ws1.Cells.Item(19, 1).Value = "03:21:41"
$ws1.Cells.Item(19, 2).Value = "04:53"
$ws1.Cells.Item(19, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(19, 4).Value = 92
$ws1.Cells.Item(19, 5).Value = "LP1912"

$ws1.Cells.Item(20, 1).Value = "03:21:41"
$ws1.Cells.Item(20, 2).Value = "05:14"
$ws1.Cells.Item(20, 3).Value = "14_ABASTO"
$ws1.Cells.Item(20, 4).Value = 113
$ws1.Cells.Item(20, 5).Value = "LP1912"

$ws1.Cells.Item(21, 1).Value = "03:21:41"
$ws1.Cells.Item(21, 2).Value = "05:16"
$ws1.Cells.Item(21, 3).Value = "17_ROMERO"
$ws1.Cells.Item(21, 4).Value = 115
$ws1.Cells.Item(21, 5).Value = "LP1912"

# ---------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Cells.Item(2, 1).Value = "Última actualización: 03:21:41"
$ws2.Cells.Item(3, 1).Value = "Total filas: 6"

$ws2.Cells.Item(10, 1).Value = "03:21:41"
$ws2.Cells.Item(10, 2).Value = "04:45"
$ws2.Cells.Item(10, 3).Value = "215A_EL PATO"
$ws2.Cells.Item(10, 4).Value = 84
$ws2.Cells.Item(10, 5).Value = "LP1912"

$ws2.Cells.Item(11, 1).Value = "02:55:01"
$ws2.Cells.Item(11, 2).Value = "04:46"
$ws2.Cells.Item(11, 3).Value = "215A_EL PATO"
$ws2.Cells.Item(11, 4).Value = 111
$ws2.Cells.Item(11, 5).Value = "LP1912"

# ---------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Cells.Item(2, 1).Value = "Última actualización: 03:21:41"
